$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 524, shifting existing rows 524:560 down to 525:561
$ws.Rows("524:524").Insert()

# Populate the newly inserted row 524 with the new record's data
$ws.Range("A524").Value = 3
$ws.Range("B524").Value = "Femacal de La Calera"
$ws.Range("C524").Value = "Coquimbo"
$ws.Range("D524").Value = 45265
$ws.Range("E524").Value = 5
$ws.Range("F524").Value = 100112001
$ws.Range("G524").Value = "Berenjena"
$ws.Range("H524").Value = "Sin especificar"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 40
$ws.Range("K524").Value = 9000
$ws.Range("L524").Value = 9000
$ws.Range("M524").Value = 9000
$ws.Range("N524").Value = "$/caja 60 unidades"
$ws.Range("O524").Value = "Región de Arica y Parinacota"
$ws.Range("P524").Value = 150
$ws.Range("Q524").Value = 60
$ws.Range("R524").Value = "Hortaliza"
